$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the "Campaña " column header: drop the trailing space.
#    (this also renames the Tabla1 table column automatically)
$ws.Range("B1").Value = "Campaña"

# 2. Add the new data row (row 4): a new agent/campaign entry.
$a4 = $ws.Cells.Item(4, 1)
$b4 = $ws.Cells.Item(4, 2)
$c4 = $ws.Cells.Item(4, 3)
$d4 = $ws.Cells.Item(4, 4)

$a4.Value = "fernandonarea6@gmail.com"
$ws.Hyperlinks.Add($a4, "mailto:fernandonarea6@gmail.com") | Out-Null
# Hyperlinks.Add() stamps its own ad-hoc formatting on the cell; restore the
# same "Hipervínculo" cell style used by the other mail cells (A2/A3).
$a4.Style = "Hipervínculo"

$b4.Value = "campaña a"
# Touch/clear wrap text so the cell keeps its own explicit style entry
# (mirrors the formatting footprint left on this cell in the source report).
$b4.WrapText = $true
$b4.WrapText = $false

$c4.Value = 654
$d4.Value = "asunto 1"

# 3. Give rows 5 and 6 an explicit (no-op) alignment touch on column D,
#    matching the extra blank, styled cells added in the source report.
$d5d6 = $ws.Range("D5:D6")
$d5d6.WrapText = $true
$d5d6.WrapText = $false

# 4. Grow "Tabla1" to include the new row.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:D4"))

# 5. Move the active selection to B1 (matches the saved cursor position).
$ws.Range("B1").Select() | Out-Null
